# edit.ps1 - Apply the "Updated cryptos list" refresh to cryptos.xlsx
#
# The sheet is a scraped cryptocurrency price/volume table (columns: A=rank,
# B=Coin, C=Link, D=Price, E=Volume(1h)). This script updates the Price and
# Volume(1h) text cells with the refreshed figures captured in the latest run,
# and reflects the two coins (rows 33/34 and 47/48) that swapped rank order.
#
# Price values are written as literal text (leading-apostrophe quote-prefix)
# so dotted/thousands-style price strings (e.g. "42.807.70") are preserved
# verbatim instead of being auto-parsed into numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '42.807.70'
$ws.Range('E2').Value = '  +0.44%  '
# Row 3
$ws.Range('D3').Value = '2.526.03'
$ws.Range('E3').Value = '  +0.42%  '
# Row 4
$ws.Range('E4').Value = '  -0.02%  '
# Row 5
$ws.Range('D5').Value = ("'" + '314.42')
$ws.Range('E5').Value = '  +1.60%  '
# Row 6
$ws.Range('D6').Value = ("'" + '95.80')
$ws.Range('E6').Value = '  -0.42%  '
# Row 7
$ws.Range('E7').Value = '  -1.76%  '
# Row 8
$ws.Range('E8').Value = '  -0.09%  '
# Row 9
$ws.Range('E9').Value = '  -0.93%  '
# Row 10
$ws.Range('D10').Value = ("'" + '36.17')
$ws.Range('E10').Value = '  -1.26%  '
# Row 11
$ws.Range('E11').Value = '  -0.38%  '
# Row 12
$ws.Range('E12').Value = '  -2.63%  '
# Row 13
$ws.Range('E13').Value = '  -3.34%  '
# Row 14
$ws.Range('D14').Value = '2.910.67'
$ws.Range('E14').Value = '  +0.35%  '
# Row 15
$ws.Range('D15').Value = ("'" + '15.28')
$ws.Range('E15').Value = '  -3.07%  '
# Row 16
$ws.Range('D16').Value = '2.479.37'
$ws.Range('E16').Value = '  -1.73%  '
# Row 17
$ws.Range('E17').Value = '  -0.44%  '
# Row 18
$ws.Range('D18').Value = '42.875.98'
$ws.Range('E18').Value = '  +0.71%  '
# Row 19
$ws.Range('D19').Value = ("'" + '12.88')
$ws.Range('E19').Value = '  -0.65%  '
# Row 20
$ws.Range('D20').Value = ("'" + '6.76')
# Row 21
$ws.Range('D21').Value = '0.0₃0964'
$ws.Range('E21').Value = '  -1.00%  '
# Row 22
$ws.Range('D22').Value = ("'" + '69.78')
$ws.Range('E22').Value = '  -2.48%  '
# Row 23
$ws.Range('D23').Value = ("'" + '254.08')
# Row 24
$ws.Range('D24').Value = ("'" + '2.96')
$ws.Range('E24').Value = '  +0.41%  '
# Row 25
$ws.Range('E25').Value = '  +1.81%  '
# Row 26
$ws.Range('D26').Value = ("'" + '26.66')
$ws.Range('E26').Value = '  -1.84%  '
# Row 27
$ws.Range('E27').Value = '  +0.09%  '
# Row 28
$ws.Range('E28').Value = '  +3.74%  '
# Row 29
$ws.Range('D29').Value = ("'" + '40.89')
$ws.Range('E29').Value = '  +8.40%  '
# Row 30
$ws.Range('D30').Value = ("'" + '10.39')
$ws.Range('E30').Value = '  +1.88%  '
# Row 31
$ws.Range('D31').Value = ("'" + '5.94')
$ws.Range('E31').Value = '  -0.16%  '
# Row 32
$ws.Range('E32').Value = '  +2.44%  '
# Row 33
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').Value = ("'" + '2.16')
$ws.Range('E33').Value = '  +3.77%  '
# Row 34
$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').Value = ("'" + '19.41')
$ws.Range('E34').Value = '  +1.36%  '
# Row 35
$ws.Range('E35').Value = '  +2.90%  '
# Row 36
$ws.Range('E36').Value = '  +0.62%  '
# Row 38
$ws.Range('E38').Value = '  -2.05%  '
# Row 39
$ws.Range('E39').Value = '  -1.23%  '
# Row 40
$ws.Range('D40').Value = ("'" + '23.59')
$ws.Range('E40').Value = '  -4.33%  '
# Row 41
$ws.Range('E41').Value = '  +14.55%  '
# Row 42
$ws.Range('E42').Value = '  +0.92%  '
# Row 43
$ws.Range('E43').Value = '  -2.18%  '
# Row 44
$ws.Range('E44').Value = '  -1.88%  '
# Row 45
$ws.Range('E45').Value = '  +0.28%  '
# Row 46
$ws.Range('D46').Value = '2.054.57'
$ws.Range('E46').Value = '  +0.92%  '
# Row 47
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').Value = ("'" + '85.55')
$ws.Range('E47').Value = '  +1.00%  '
# Row 48
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = ("'" + '108.81')
$ws.Range('E48').Value = '  +7.02%  '
# Row 49
$ws.Range('E49').Value = '  -0.67%  '
# Row 50
$ws.Range('D50').Value = ("'" + '75.46')
$ws.Range('E50').Value = '  +3.59%  '
# Row 51
$ws.Range('D51').Value = '2.764.92'
$ws.Range('E51').Value = '  +0.27%  '
